$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) / Volume(1h) (E) figures for the crypto list refresh.
$updates = @(
    @{ Row = 2; D = '27.966.98'; E = '  -0.22%  ' },
    @{ Row = 3; D = '1.859.65'; E = '  -0.58%  ' },
    @{ Row = 4; D = '1.004'; E = '  -0.06%  ' },
    @{ Row = 5; D = '311.96'; E = $null },
    @{ Row = 6; D = '1.002'; E = '  -0.06%  ' },
    @{ Row = 7; D = '0.5139'; E = '  +1.96%  ' },
    @{ Row = 8; D = '0.3817'; E = '  -0.51%  ' },
    @{ Row = 9; D = '0.08248'; E = '  -4.41%  ' },
    @{ Row = 10; D = $null; E = '  -0.40%  ' },
    @{ Row = 11; D = '41.49'; E = '  -0.09%  ' },
    @{ Row = 12; D = '6.183'; E = '  -1.92%  ' },
    @{ Row = 13; D = '20.47'; E = '  -0.72%  ' },
    @{ Row = 14; D = '1.852.15'; E = '  -1.91%  ' },
    @{ Row = 15; D = '7.279'; E = '  +1.54%  ' },
    @{ Row = 16; D = $null; E = '  -0.07%  ' },
    @{ Row = 17; D = '0.00001093'; E = '  -0.48%  ' },
    @{ Row = 18; D = '90.31'; E = '  -0.62%  ' },
    @{ Row = 19; D = '0.06640'; E = '  +0.21%  ' },
    @{ Row = 20; D = '17.66'; E = '  -1.70%  ' },
    @{ Row = 22; D = '6.014'; E = '  -0.93%  ' },
    @{ Row = 23; D = '28.003.05'; E = '  -0.25%  ' },
    @{ Row = 24; D = $null; E = '  -2.89%  ' },
    @{ Row = 25; D = '2.243'; E = '  -0.89%  ' },
    @{ Row = 26; D = '2.069.43'; E = '  -1.19%  ' },
    @{ Row = 27; D = '2.503'; E = '  -2.15%  ' },
    @{ Row = 28; D = '157.29'; E = '  +0.26%  ' },
    @{ Row = 29; D = $null; E = '  -1.16%  ' },
    @{ Row = 30; D = '124.37'; E = '  -1.21%  ' },
    @{ Row = 31; D = $null; E = '  +1.08%  ' },
    @{ Row = 32; D = $null; E = '  -2.97%  ' },
    @{ Row = 33; D = '5.824'; E = '  +4.13%  ' },
    @{ Row = 34; D = '3.594'; E = '  -0.13%  ' },
    @{ Row = 35; D = '9.386'; E = '  -2.30%  ' },
    @{ Row = 36; D = '0.02407'; E = '  -0.94%  ' },
    @{ Row = 37; D = '0.06486'; E = '  -1.09%  ' },
    @{ Row = 38; D = '0.2191'; E = '  +1.16%  ' },
    @{ Row = 39; D = '0.6533'; E = '  +2.93%  ' },
    @{ Row = 40; D = $null; E = '  -0.82%  ' },
    @{ Row = 41; D = '4.980'; E = '  +2.30%  ' },
    @{ Row = 42; D = '1.207'; E = '  -2.28%  ' },
    @{ Row = 43; D = $null; E = '  -3.05%  ' },
    @{ Row = 44; D = '0.6105'; E = '  +2.23%  ' },
    @{ Row = 45; D = '13.00'; E = '  -1.02%  ' },
    @{ Row = 46; D = $null; E = '  -0.02%  ' },
    @{ Row = 47; D = $null; E = '  -0.44%  ' },
    @{ Row = 48; D = '2.011'; E = '  +1.45%  ' },
    @{ Row = 49; D = $null; E = '  -1.38%  ' },
    @{ Row = 50; D = '120.66'; E = '  -0.51%  ' },
    @{ Row = 51; D = '77.73'; E = '  -3.16%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Force text storage so values like "1.004" or "27.966.98" keep their
        # original (non-numeric) literal formatting instead of becoming numbers.
        $ws.Range("D" + $u.Row).NumberFormat = "@"
        $ws.Range("D" + $u.Row).Value = $u.D
        $ws.Range("D" + $u.Row).Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
